$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("n5", "n5_IMG_3176.jpeg", "True", "no_meltpatch", "negative"),
    @("n6", "n6_IMG_3180.jpeg", "True", "no_meltpatch", "negative"),
    @("n7", "n7_IMG_3179.jpeg", "True", "no_meltpatch", "negative"),
    @("n8", "n8_IMG_3175.jpeg", "True", "no_meltpatch", "negative")
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        $value = $rowData[$c]
        # Values like "True"/"False" would otherwise be auto-coerced to a
        # boolean by COM's type inference. Prefix with an apostrophe to force
        # literal text entry (matches the source data, which stores "True"
        # as a plain string), then strip the resulting quote-prefix style so
        # the cell format matches the other text cells in the column.
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    }
}
